$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Range("E1").Value = "Date Created (Year)*"
$ws.Range("E2").Value = 2000
$ws.Range("E1:E4").Font.Color = 0
$ws.Range("E3:E4").Select()
